$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")
$ws.Cells.Item(84, 1).NumberFormat = "@"
$ws.Cells.Item(84, 1).Value = "2025-12-27"
$ws.Cells.Item(84, 2).Value = 0.0
$ws.Cells.Item(84, 3).Value = 28.0
